$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Customer #1 (row 1) ---------------------------------------------------
# A new "NaamKlant0001" (customer-name) header replaces the old L1 value and
# the numeric code in F1 gets an extra two trailing zeros (11 -> 1100).
$ws.Range("F1").Value = 1100
$ws.Range("L1").Value = "NaamKlant0001"

# --- Customer #2 (row 2) ---------------------------------------------------
# Same pair of edits for the second record.
$ws.Range("F2").Value = 2200
$ws.Range("L2").Value = "NaamKlant0002"

# --- View / layout ----------------------------------------------------------
# Scroll the sheet so column F is the left-most visible column (topLeftCell
# moves from A1 to F1) and widen the data columns (A:AE) so the new, longer
# labels are readable; columns beyond AE keep the original default width.
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("A1:AE2").ColumnWidth = 19.5
